$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "F003"
$ws.Range("B4").Value = "??"
$ws.Range("C4").Value = "Dalock/Bauß"
$ws.Range("D4").Value = "Stellen/Transitionen lassen sich in der GUI (Im Regel-Fenster) untereinander schieben"
$ws.Range("E4").Value = "im Verzeichnis F003_F004"
$ws.Range("F4").Value = 2

$ws.Range("A5").Value = "F004"
$ws.Range("B5").Value = "??"
$ws.Range("C5").Value = "Dalock/Bauß"
$ws.Range("D5").Value = "Beim Löschen eines Elementes aus der GUI wird, wenn man danach einen Rechtsklick ausführt eine oder mehrere Exceptions geworfen"
$ws.Range("E5").Value = "im Verzeichnis F003_F005"
$ws.Range("F5").Value = 2

$ws.Range("D6").Select()
